$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-code" (column D) and "codeforiati:group-name"
# (column E) columns were swapped throughout the data range (including
# the header row), i.e. D<->E content exchange for rows 1-94.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 1) { $lastRow = 94 }

for ($r = 1; $r -le 94; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
